$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I33").Value = 129.9
$ws.Range("J33").Value = 1000000000
$ws.Range("K33").Value = 129.9
$ws.Range("L33").Value = 1000000000
$ws.Range("M33").Value = 99.09999999999999
$ws.Range("N33").Value = -1000000458

$ws.Range("H40").Value = 2483
$ws.Range("I40").Value = 2355.375
$ws.Range("J40").Value = 2993.5
$ws.Range("K40").Value = 2355.375
$ws.Range("L40").Value = 2993.5
$ws.Range("M40").Value = -2180.375
$ws.Range("N40").Value = -3343.5

$ws.Range("H43").Value = 11405.4
$ws.Range("I43").Value = 12095.667
$ws.Range("J43").Value = 10370
$ws.Range("K43").Value = 12095.667
$ws.Range("L43").Value = 10370
$ws.Range("M43").Value = -12026.667
$ws.Range("N43").Value = -10508

$ws.Range("H45").Value = 2833.3333
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 2833.3333
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 8499.999899999999
$ws.Range("M45").Value = ""
$ws.Range("N45").Value = -8883.999899999999

$ws.Range("H52").Value = 1974.1
$ws.Range("I52").Value = 534.7143
$ws.Range("J52").Value = 5332.6665
$ws.Range("K52").Value = 1604.1429
$ws.Range("L52").Value = 15997.9995
$ws.Range("M52").Value = -1444.1429
$ws.Range("N52").Value = -16317.9995

$ws.Range("H137").Value = 2518.375
$ws.Range("I137").Value = 2552.9333
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 7658.7999
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -5108.7999
$ws.Range("N137").Value = -11100

$ws.Range("H138").Value = 1715.9344
$ws.Range("I138").Value = 818.2
$ws.Range("J138").Value = 2339.361
$ws.Range("K138").Value = 2454.6
$ws.Range("L138").Value = 7018.083
$ws.Range("M138").Value = 2685.4
$ws.Range("N138").Value = -17298.083

$ws.Range("H141").Value = 6012.0835
$ws.Range("I141").Value = 4614.7
$ws.Range("J141").Value = 12999
$ws.Range("K141").Value = 13844.1
$ws.Range("L141").Value = 38997
$ws.Range("M141").Value = -8664.099999999999
$ws.Range("N141").Value = -49357

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1932.8
$ws.Range("I61").Value = 1932.8
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1932.8
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1720.8
$ws.Range("N61").Value = ""

$ws.Range("H63").Value = 4487
$ws.Range("I63").Value = 4932.3335
$ws.Range("J63").Value = 3151
$ws.Range("K63").Value = 4932.3335
$ws.Range("L63").Value = 3151
$ws.Range("M63").Value = -4246.3335

$ws.Range("H66").Value = 4487
$ws.Range("I66").Value = 4932.3335
$ws.Range("J66").Value = 3151
$ws.Range("K66").Value = 24661.6675
$ws.Range("L66").Value = 15755
$ws.Range("M66").Value = -21229.6675

$ws.Range("H74").Value = 3056.1
$ws.Range("I74").Value = 3074.8
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 3074.8
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -2200.8
$ws.Range("N74").Value = -4748

$ws.Range("H77").Value = 3056.1
$ws.Range("I77").Value = 3074.8
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 15374
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -11006
$ws.Range("N77").Value = -23736

$ws.Range("H97").Value = 464.42307
$ws.Range("I97").Value = 267.42105
$ws.Range("J97").Value = 999.1429000000001
$ws.Range("K97").Value = 267.42105
$ws.Range("L97").Value = 999.1429000000001
$ws.Range("M97").Value = 228.57895
$ws.Range("N97").Value = -1991.1429

$ws.Range("H102").Value = 22737.666
$ws.Range("I102").Value = 6390.231
$ws.Range("J102").Value = 128996
$ws.Range("K102").Value = 6390.231
$ws.Range("L102").Value = 128996
$ws.Range("M102").Value = -4768.231
$ws.Range("N102").Value = -132240

$ws.Range("H132").Value = 7338.825
$ws.Range("I132").Value = 4724.737
$ws.Range("J132").Value = 57006.5
$ws.Range("K132").Value = 14174.211
$ws.Range("L132").Value = 171019.5
$ws.Range("M132").Value = -11644.211

$ws.Range("H136").Value = 1932.8
$ws.Range("I136").Value = 1932.8
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5798.4
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3248.4
$ws.Range("N136").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1954.7333
$ws.Range("I105").Value = 1682.65
$ws.Range("J105").Value = 2498.9
$ws.Range("K105").Value = 1682.65
$ws.Range("L105").Value = 2498.9
$ws.Range("M105").Value = 64.34999999999991

$ws.Range("H134").Value = 3882.5833
$ws.Range("I134").Value = 3842.0952
$ws.Range("J134").Value = 4166
$ws.Range("K134").Value = 11526.2856
$ws.Range("L134").Value = 12498
$ws.Range("M134").Value = -8991.285600000001
$ws.Range("N134").Value = -17568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 40.3125
$ws.Range("I7").Value = 28.666666
$ws.Range("J7").Value = 55.285713
$ws.Range("K7").Value = 28.666666
$ws.Range("L7").Value = 55.285713
$ws.Range("M7").Value = 84.33333400000001
$ws.Range("N7").Value = -281.285713

$ws.Range("H58").Value = 3021
$ws.Range("I58").Value = 3468.1667
$ws.Range("J58").Value = 2637.7144
$ws.Range("K58").Value = 3468.1667
$ws.Range("L58").Value = 2637.7144
$ws.Range("M58").Value = -3265.1667

$ws.Range("H132").Value = 3485.875
$ws.Range("I132").Value = 3195.6
$ws.Range("J132").Value = 3969.6667
$ws.Range("K132").Value = 9586.799999999999
$ws.Range("L132").Value = 11909.0001
$ws.Range("M132").Value = -7056.799999999999
$ws.Range("N132").Value = -16969.0001

$ws.Range("H134").Value = 2434.8333
$ws.Range("I134").Value = 2152.25
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 6456.75
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -3921.75

$ws.Range("H136").Value = 3021
$ws.Range("I136").Value = 3468.1667
$ws.Range("J136").Value = 2637.7144
$ws.Range("K136").Value = 10404.5001
$ws.Range("L136").Value = 7913.1432
$ws.Range("M136").Value = -7854.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 335032.66
$ws.Range("I32").Value = 501299.5
$ws.Range("J32").Value = 2499
$ws.Range("K32").Value = 1503898.5
$ws.Range("L32").Value = 7497
$ws.Range("M32").Value = -1503615.5

$ws.Range("H54").Value = 7166
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 7166
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 21498
$ws.Range("N54").Value = -22616

$ws.Range("H69").Value = 6100
$ws.Range("I69").Value = 6675
$ws.Range("J69").Value = 1500
$ws.Range("K69").Value = 20025
$ws.Range("L69").Value = 4500
$ws.Range("M69").Value = -19214

$ws.Range("H72").Value = 6100
$ws.Range("I72").Value = 6675
$ws.Range("J72").Value = 1500
$ws.Range("K72").Value = 60075
$ws.Range("L72").Value = 13500
$ws.Range("M72").Value = -56019

$ws.Range("H131").Value = 1290.6666
$ws.Range("I131").Value = 828.2727
$ws.Range("J131").Value = 1799.3
$ws.Range("K131").Value = 2484.8181
$ws.Range("L131").Value = 5397.9
$ws.Range("M131").Value = 2555.1819
$ws.Range("N131").Value = -15477.9

$ws.Range("H138").Value = 4062.125
$ws.Range("I138").Value = 4062.125
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 12186.375
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -7046.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 410.95
$ws.Range("I107").Value = 298.7647
$ws.Range("J107").Value = 1046.6666
$ws.Range("K107").Value = 896.2941000000001
$ws.Range("L107").Value = 3139.9998
$ws.Range("M107").Value = 1023.7059
$ws.Range("N107").Value = -6979.9998

$ws.Range("H136").Value = 5109.2104
$ws.Range("I136").Value = 2948.9375
$ws.Range("J136").Value = 16630.666
$ws.Range("K136").Value = 8846.8125
$ws.Range("L136").Value = 49891.99800000001
$ws.Range("M136").Value = -6296.8125
$ws.Range("N136").Value = -54991.99800000001
